$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray C2 value (y_0_forecast for 2008 row) - cell becomes empty
$ws.Range("C2").ClearContents()

# Corrected y_0_forecast (column C) values - tiny floating point fixes from the
# naive component forecaster bug fix
$ws.Range("C3").Value = -4.857241224140941
$ws.Range("C4").Value = 0.3625742673738719
$ws.Range("C6").Value = 0.6836026627130787
$ws.Range("C7").Value = 0.2336505480021511
$ws.Range("C9").Value = -0.1588690085688071
$ws.Range("C10").Value = -0.5438176183081955
$ws.Range("C11").Value = -0.006876704825686808
$ws.Range("C13").Value = -0.63478973259814
$ws.Range("C14").Value = -0.7158018152081613
$ws.Range("C15").Value = 1.576357831383679
$ws.Range("C16").Value = -1.488707312182613
$ws.Range("C18").Value = 0.8934739937295433
$ws.Range("C19").Value = -1.372720900450863

# Corrected y_1_forecast (column E) values
$ws.Range("E8").Value = -0.03923323971221082
$ws.Range("E9").Value = -0.4617525814883061
$ws.Range("E10").Value = 0.01247916696662799
$ws.Range("E13").Value = 0.01241557525979431
$ws.Range("E14").Value = -0.7615805088034833
$ws.Range("E16").Value = -1.062239424572287
$ws.Range("E18").Value = 0.4129745242491101
